$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = [double]"1"
$ws.Range("E2").Value = [double]"1"

# Row 3
$ws.Range("D3").Value = [double]"2.50121850521859E-06"
$ws.Range("E3").Value = [double]"2.50121850521859E-06"

# Row 4
$ws.Range("D4").Value = [double]"7.826577561566211E-06"
$ws.Range("E4").Value = [double]"7.826577561566211E-06"

# Row 5
$ws.Range("D5").Value = [double]"0.03953455114500491"
$ws.Range("E5").Value = [double]"0.03953455114500491"

# Row 6
$ws.Range("D6").Value = [double]"2.615731264744623E-22"
$ws.Range("E6").Value = [double]"2.615731264744623E-22"

# Row 7
$ws.Range("D7").Value = [double]"0.9999999999443041"
$ws.Range("E7").Value = [double]"5.569589234255545E-11"

# Row 8
$ws.Range("D8").Value = [double]"1.02496365851681E-09"
$ws.Range("E8").Value = [double]"0.9999999989750363"

# Row 9
$ws.Range("D9").Value = [double]"1.164990038683197E-05"
$ws.Range("E9").Value = [double]"0.9999883500996132"

# Row 10
$ws.Range("D10").Value = [double]"0.9999820694967964"
$ws.Range("E10").Value = [double]"1.793050320364475E-05"

# Row 11
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = [double]"0.9873128950275537"
$ws.Range("E11").Value = [double]"0.01268710497244629"
$ws.Range("F11").Value = [double]"8.150106430053711"
$ws.Range("G11").Value = [double]"0.7"
